$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Ibuprofen's Initial Stock (B3) from 50 to 9
$ws.Range("B3").Value = 9

# Move the active cell/selection to B3
$ws.Range("B3").Select()
